$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.983.17"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "3.384.91"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.54"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.392"
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").Value = "3.964.54"
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.126"
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.62"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("E15").Value = "  -1.25%  "
$ws.Range("D16").Value = "3.335.78"
$ws.Range("E16").Value = "  -2.98%  "
$ws.Range("D17").Value = "61.057.86"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.10"
$ws.Range("E18").Value = "  -4.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.71"
$ws.Range("E19").Value = "  -5.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.91"
$ws.Range("E20").Value = "  -5.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.04"
$ws.Range("E21").Value = "  -4.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.76"
$ws.Range("E23").Value = "  -2.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  -4.56%  "
$ws.Range("D26").Value = "3.528.77"
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.180"
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -2.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.15"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.96"
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.41"
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.40"
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.95"
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.41"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").Value = "3.417.26"
$ws.Range("E37").Value = "  -1.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.47"
$ws.Range("E39").Value = "  -5.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0768"
$ws.Range("E40").Value = "  -2.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.24"
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.779"
$ws.Range("E43").Value = "  -2.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.82"
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("E45").Value = "  -2.57%  "
$ws.Range("E46").Value = "  -3.25%  "
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("D48").Value = "2.457.82"
$ws.Range("E48").Value = "  -6.14%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.76"
$ws.Range("E49").Value = "  -2.95%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.90"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("E51").Value = "  +2.50%  "
